$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.797.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.726.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.13%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '504.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  -3.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.734.13'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.348'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.126'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.207.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.908.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.75%  '
$ws.Range("E16").Value = '  -4.18%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.27%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.724.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.26%  '
$ws.Range("E19").Value = '  -3.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '344.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.41%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.13'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -4.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.171'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("E29").Value = '  -3.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0832'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  -1.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.18'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.951'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.97%  '
$ws.Range("E40").Value = '  -6.26%  '
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.187.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.995'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0555'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.601'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.58%  '
$ws.Range("E46").Value = '  -7.94%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0225'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0885'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.82%  '
